$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list (price/volume refresh + two rank swaps: Chainlink/Polygon, Maker/Stacks)
# Column D ("Price") values are stored as text in this sheet (dotted/locale formatted
# numbers, e.g. "46.528.59"), so new values are written with a leading apostrophe to
# force Excel to keep them as text instead of auto-converting to a number.

# Row 2
$ws.Range("D2").Value = "'46.488.79"
$ws.Range("E2").Value = "  +5.32%  "
# Row 3
$ws.Range("D3").Value = "'2.302.48"
$ws.Range("E3").Value = "  +4.21%  "
# Row 4
$ws.Range("E4").Value = "  -0.58%  "
# Row 5
$ws.Range("D5").Value = "'299.43"
$ws.Range("E5").Value = "  +0.12%  "
# Row 6
$ws.Range("D6").Value = "'97.09"
$ws.Range("E6").Value = "  +7.70%  "
# Row 7
$ws.Range("D7").Value = "'0.573"
$ws.Range("E7").Value = "  +0.19%  "
# Row 8
$ws.Range("E8").Value = "  -0.46%  "
# Row 9
$ws.Range("D9").Value = "'0.528"
$ws.Range("E9").Value = "  +7.96%  "
# Row 10
$ws.Range("D10").Value = "'35.53"
$ws.Range("E10").Value = "  +6.40%  "
# Row 11
$ws.Range("D11").Value = "'0.0800"
$ws.Range("E11").Value = "  +1.97%  "
# Row 12
$ws.Range("D12").Value = "'7.36"
$ws.Range("E12").Value = "  +7.58%  "
# Row 13
$ws.Range("E13").Value = "  +0.86%  "
# Row 14
$ws.Range("D14").Value = "'2.654.33"
$ws.Range("E14").Value = "  +4.30%  "
# Row 15
$ws.Range("D15").Value = "'2.301.53"
$ws.Range("E15").Value = "  +3.94%  "
# Row 16
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'13.96"
$ws.Range("E16").Value = "  +4.49%  "
# Row 17
$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").Value = "'0.824"
$ws.Range("E17").Value = "  +5.82%  "
# Row 18
$ws.Range("D18").Value = "'46.460.92"
$ws.Range("E18").Value = "  +5.40%  "
# Row 19
$ws.Range("D19").Value = "'13.13"
$ws.Range("E19").Value = "  +18.67%  "
# Row 20
$ws.Range("D20").Value = "'0.0₃0940"
$ws.Range("E20").Value = "  +5.14%  "
# Row 21
$ws.Range("D21").Value = "'6.14"
$ws.Range("E21").Value = "  +4.15%  "
# Row 22
$ws.Range("D22").Value = "'66.90"
$ws.Range("E22").Value = "  +4.78%  "
# Row 23
$ws.Range("D23").Value = "'247.97"
$ws.Range("E23").Value = "  +7.21%  "
# Row 24
$ws.Range("D24").Value = "'2.92"
$ws.Range("E24").Value = "  +3.59%  "
# Row 25
$ws.Range("D25").Value = "'1.99"
$ws.Range("E25").Value = "  +7.00%  "
# Row 26
$ws.Range("E26").Value = "  -0.22%  "
# Row 27
$ws.Range("D27").Value = "'42.71"
$ws.Range("E27").Value = "  +17.08%  "
# Row 28
$ws.Range("D28").Value = "'2.28"
$ws.Range("E28").Value = "  +2.21%  "
# Row 29
$ws.Range("D29").Value = "'9.81"
$ws.Range("E29").Value = "  +5.58%  "
# Row 30
$ws.Range("D30").Value = "'20.07"
$ws.Range("E30").Value = "  +4.43%  "
# Row 31
$ws.Range("D31").Value = "'5.74"
$ws.Range("E31").Value = "  +5.91%  "
# Row 32
$ws.Range("D32").Value = "'0.0799"
$ws.Range("E32").Value = "  +7.29%  "
# Row 33
$ws.Range("D33").Value = "'146.23"
$ws.Range("E33").Value = "  +0.55%  "
# Row 34
$ws.Range("E34").Value = "  +3.26%  "
# Row 35
$ws.Range("E35").Value = "  +6.87%  "
# Row 36
$ws.Range("D36").Value = "'0.111"
$ws.Range("E36").Value = "  +7.34%  "
# Row 37
$ws.Range("E37").Value = "  +2.10%  "
# Row 38
$ws.Range("D38").Value = "'1.80"
$ws.Range("E38").Value = "  +8.94%  "
# Row 39
$ws.Range("D39").Value = "'3.99"
$ws.Range("E39").Value = "  +12.84%  "
# Row 40
$ws.Range("D40").Value = "'14.71"
$ws.Range("E40").Value = "  +7.17%  "
# Row 41
$ws.Range("E41").Value = "  +7.69%  "
# Row 42
$ws.Range("D42").Value = "'0.0306"
$ws.Range("E42").Value = "  +7.78%  "
# Row 43
$ws.Range("D43").Value = "'0.998"
$ws.Range("E43").Value = "  -0.89%  "
# Row 44
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "'1.840.27"
$ws.Range("E44").Value = "  +4.17%  "
# Row 45
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").Value = "'1.96"
$ws.Range("E45").Value = "  +18.49%  "
# Row 46
$ws.Range("D46").Value = "'89.57"
$ws.Range("E46").Value = "  +19.69%  "
# Row 47
$ws.Range("D47").Value = "'0.196"
$ws.Range("E47").Value = "  +12.41%  "
# Row 48
$ws.Range("D48").Value = "'72.13"
$ws.Range("E48").Value = "  +3.85%  "
# Row 49
$ws.Range("D49").Value = "'4.93"
$ws.Range("E49").Value = "  +11.18%  "
# Row 50
$ws.Range("D50").Value = "'97.22"
# Row 51
$ws.Range("D51").Value = "'54.38"
$ws.Range("E51").Value = "  +8.08%  "
